# Commit: "rules need new pdf"
# Fixes two typos in the card-effect descriptions and resets the
# sheet's scroll/selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A7 holds the "Redraw" card text - add a trailing period
$ws.Range("A7").Value = "Redraw: Return any amount of tiles to the bag and draw new ones."

# A1 holds the "Bomb" card text - fix "adjecent" -> "adjacent"
$ws.Range("A1").Value = "Bomb: Destroy a tile and its adjacent ones (+ shape).  Tiles go back to pile. Can't destroy tile under player.  "

# Reset the view: move the selection to A6 (previously the sheet was
# scrolled to A7 with A17 selected).
$ws.Range("A6").Select()
